# Generate Report for Handback
# - The "78beadd6..." row (row 3) status flips from "Ready for handoff" to
#   "Handback transform failed" on the Overview sheet and on each per-locale
#   sheet (the same shared text is shown in Overview!E3/F3 and in
#   zh-cn!C3 / de-de!C3).
# - Each per-locale sheet gets a new "Error Detail" message in column P,
#   row 3, describing the handback/handoff filename mismatch.
# - The "Error Detail" column (P) is widened to fit the new message text.

$wb = $excel.ActiveWorkbook

$statusText = "Handback transform failed"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

# Excel's ColumnWidth property (character units) round-trips through the
# stored <col width=".."/> (Maximum Digit Width units) with a constant
# +5/6 padding offset, so request 5/6 less to land exactly on 40 in the
# saved file.
$targetColWidth = 40 - (5 / 6)

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $statusText
$wsZhCn.Range("P3").Value = "Handback file name: hn5a0pmg.jtl is different with handoff file name: 78beadd6-f4e2-43f9-80c8-6e4a888ea7ba.e546683f3b69eb90fc681d5d4dd1aa96583d6a9f.zh-cn."
$wsZhCn.Columns.Item(16).ColumnWidth = $targetColWidth

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $statusText
$wsDeDe.Range("P3").Value = "Handback file name: hn5a0pmg.jtl is different with handoff file name: 78beadd6-f4e2-43f9-80c8-6e4a888ea7ba.e546683f3b69eb90fc681d5d4dd1aa96583d6a9f.de-de."
$wsDeDe.Columns.Item(16).ColumnWidth = $targetColWidth
